# Automatic update of files.
# 1) Bump the "Förändrad" (Changed) date in column C from 45184 to 45186
#    for every data row (rows 2..205).
# 2) For the rows whose HYPERLINK formulas still have only one argument
#    (rows 2..4, columns S, T, V, W, X, Y), add the friendly display
#    name (the value in column A of that row) as the second HYPERLINK
#    argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count - $usedRange.Row + 1
if ($lastRow -lt 205) { $lastRow = 205 }

$oldDate = 45184
$newDate = 45186

$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq $oldDate) {
        $cCell.Value = $newDate
    }

    $name = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$r")
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) { continue }
        if ($formula.IndexOf("HYPERLINK(") -lt 0) { continue }
        if ($formula.IndexOf(",") -ge 0) { continue }

        $trimmed = $formula.TrimEnd()
        if (-not $trimmed.EndsWith(')')) { continue }

        $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $name + '")'
        $cell.Formula = $newFormula
    }
}
